$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J49
$iVals = @(9,8,7,7,9,9,8,8,9,9,6,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,5,7,3,5,3)
$jVals = @(9,9,7,8,9,9,8,8,9,9,6,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,5,7,3,5,3)

for ($r = 0; $r -lt 48; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
